# CAD-1153: include the filters and templates to show the 3YC assets and request
# Adds three new trailing columns (commitment / commitment start date / commitment
# end date) to the "Data" sheet header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New header cells right after the existing last column (AI).
$ws.Range("AJ1").Value = "commitment"
$ws.Range("AK1").Value = "commitment start date"
$ws.Range("AL1").Value = "commitment end date"

# Size the 3 new columns to fit their header text (mirrors the bestFit columns
# already present on the sheet, e.g. L/M/Q..AH).
$ws.Columns.Item(36).ColumnWidth = 10.330729166666666
$ws.Columns.Item(37).ColumnWidth = 17.998697916666668
$ws.Columns.Item(38).ColumnWidth = 17.330729166666668

# Leave the selection where the author ended up after adding the columns.
$ws.Range("AG22").Select()
